# Update computed values in result_data_RandomForest.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -7.033499999999992
$ws.Range("C12").Value = -11.38799999999999
$ws.Range("D14").Value = -7.722500000000004
$ws.Range("D26").Value = -8.4405
$ws.Range("D31").Value = -8.467199999999998
$ws.Range("C32").Value = -13.37090000000001
$ws.Range("D35").Value = -8.208299999999998
$ws.Range("C36").Value = -12.7968
$ws.Range("D37").Value = -7.811400000000003
$ws.Range("C38").Value = -12.5401
$ws.Range("D45").Value = -7.813599999999997
$ws.Range("C46").Value = -14.6781
$ws.Range("C54").Value = -12.6347
$ws.Range("C55").Value = -13.8209
$ws.Range("D57").Value = -8.3331
$ws.Range("C67").Value = -11.0797
$ws.Range("C69").Value = -12.22209999999999
$ws.Range("C72").Value = -11.47410000000001
$ws.Range("C91").Value = -10.5457
$ws.Range("C99").Value = -12.97499999999999
$ws.Range("D100").Value = -7.981200000000005
$ws.Range("D102").Value = -7.826400000000001
